# Weekly update: a new week of Mango price data is added at the top of the
# data block (row 48), pushing all existing data rows down by one. The row
# that previously fell off the bottom of the sheet (old row 157) becomes the
# new last row (158).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 48 - this shifts rows 48:157 down to 49:158
# (carrying their formatting/values with them) and extends the used range.
$ws.Rows("48:48").Insert()

# Populate the newly inserted row 48 with this week's record.
$ws.Cells.Item(48, 1).Value  = 11
$ws.Cells.Item(48, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(48, 3).Value  = "Bíobío"
$ws.Cells.Item(48, 4).Value  = 45028
$ws.Cells.Item(48, 5).Value  = 8
$ws.Cells.Item(48, 6).Value  = "Fruta"
$ws.Cells.Item(48, 7).Value  = 100108
$ws.Cells.Item(48, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(48, 9).Value  = 100108002
$ws.Cells.Item(48, 10).Value = "Mango"
$ws.Cells.Item(48, 11).Value = "Sin especificar"
$ws.Cells.Item(48, 12).Value = "Primera"
$ws.Cells.Item(48, 13).Value = 200
$ws.Cells.Item(48, 14).Value = 7000
$ws.Cells.Item(48, 15).Value = 7500
$ws.Cells.Item(48, 16).Value = 7250
$ws.Cells.Item(48, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(48, 18).Value = "Perú"
$ws.Cells.Item(48, 19).Value = 1812
$ws.Cells.Item(48, 20).Value = 4
